# Update the CDA Logical model (ValueSet-CDAActClass) for ST.r2b
$wb = $excel.ActiveWorkbook

# 1. Rename the include sheet
$includeSheet = $wb.Worksheets.Item("Include from ActClass")
$includeSheet.Name = "Include #0"

# 2. Update the Metadata sheet
$ws = $wb.Worksheets.Item("Metadata")

# Bump the Version and Date property values
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" property row right after "Contact" (row 10)
# and before "Description" (old row 11) - matching formatting of the row above.
$ws.Rows.Item(11).Insert()
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").ClearContents()
